$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Riley Meredith"

# Insert a new column before column A, shifting existing data right
$ws.Range("A1").EntireColumn.Insert()

# Fill in the new column header and value
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "8th"
